$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.069.85'
$ws.Range("E2").Value = '  -0.55%  '

$ws.Range("D3").Value = '2.613.90'
$ws.Range("E3").Value = '  -1.08%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.91'
$ws.Range("E5").Value = '  -1.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.24'
$ws.Range("E6").Value = '  -1.16%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.531'
$ws.Range("E8").Value = '  -2.44%  '

$ws.Range("D9").Value = '2.613.42'
$ws.Range("E9").Value = '  -1.09%  '

$ws.Range("E10").Value = '  -4.61%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.161'
$ws.Range("E11").Value = '  +1.58%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.362'
$ws.Range("E12").Value = '  -0.80%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.19'
$ws.Range("E13").Value = '  -0.72%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.28'
$ws.Range("E14").Value = '  -2.70%  '

$ws.Range("D15").Value = '3.091.40'
$ws.Range("E15").Value = '  -1.05%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000179'
$ws.Range("E16").Value = '  -2.73%  '

$ws.Range("D17").Value = '67.015.85'
$ws.Range("E17").Value = '  -0.45%  '

$ws.Range("D18").Value = '2.609.05'
$ws.Range("E18").Value = '  -1.29%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.74'
$ws.Range("E19").Value = '  -1.57%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.78'
$ws.Range("E20").Value = '  -0.94%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '354.26'
$ws.Range("E21").Value = '  -2.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.27'
$ws.Range("E22").Value = '  -3.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.62'
$ws.Range("E23").Value = '  -3.49%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.45'
$ws.Range("E24").Value = '  -4.93%  '

$ws.Range("E25").Value = '  +0.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.91'
$ws.Range("E26").Value = '  -5.23%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '69.06'
$ws.Range("E27").Value = '  -2.76%  '

$ws.Range("E28").Value = '  -1.17%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("D30").Value = '0.0₃0993'
$ws.Range("E30").Value = '  -3.31%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '542.57'
$ws.Range("E31").Value = '  -2.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.84'
$ws.Range("E32").Value = '  -3.01%  '

$ws.Range("E33").Value = '  -3.76%  '

$ws.Range("E34").Value = '  -3.03%  '

$ws.Range("E35").Value = '  +0.85%  '

$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.48'
$ws.Range("E37").Value = '  -3.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '158.99'
$ws.Range("E38").Value = '  +0.79%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.91'
$ws.Range("E39").Value = '  -2.67%  '

$ws.Range("E40").Value = '  -2.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.24'
$ws.Range("E41").Value = '  +1.72%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.80'
$ws.Range("E42").Value = '  -1.50%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.13'
$ws.Range("E43").Value = '  -2.66%  '

$ws.Range("E44").Value = '  +0.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.40'
$ws.Range("E45").Value = '  -4.90%  '

$ws.Range("D46").Value = '0.0₆0298'
$ws.Range("E46").Value = '  -1.03%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '151.26'
$ws.Range("E47").Value = '  -1.71%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.574'
$ws.Range("E48").Value = '  -3.70%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.76'
$ws.Range("E49").Value = '  -3.36%  '

$ws.Range("E50").Value = '  -1.81%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0768'
$ws.Range("E51").Value = '  -1.36%  '
